# Generate Report for Handoff
# - Updates the "Status" text from "Handed back: in sync with en-US" to
#   "Ready for handoff" on all three sheets (Overview, zh-cn, de-de).
# - Refreshes the related "Latest HO Xliff Generate Date" / "Latest Handoff
#   Datetime" timestamps to reflect the new handoff run.
# - Narrows the now-shorter Status/date columns to fit the shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn (E2) and de-de (F2) status cells, plus the
#     shared "Latest HO Xliff Generate Date" (G2). ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-27 19:06:14"

# --- zh-cn detail sheet: Status (C2) and Latest Handoff Datetime (H2). ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-27 19:06:10"

# --- de-de detail sheet: Status (C2) and Latest Handoff Datetime (H2). ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-27 19:06:14"

# --- Narrow the Status/date columns that held the old, longer text.
#     (ColumnWidth is expressed in characters; 16.3 is the value that maps
#     to the target stored column width of ~17.22 used by the handoff
#     report template.) ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
